$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the marking value per correct answer (row 11, "Marking")
$ws.Range("B11").Value = 5

# Update the total score (row 12, "Total")
$ws.Range("B12").Value = 35

# Update the correct/total marks summary text
$ws.Range("E12").Value = "35/140"
